$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1642.1428
$ws.Range("I111").Value = 1039
$ws.Range("J111").Value = 3150
$ws.Range("K111").Value = 3117
$ws.Range("L111").Value = 9450
$ws.Range("M111").Value = -50
$ws.Range("N111").Value = -15584

$ws.Range("H127").Value = 1273.4546
$ws.Range("I127").Value = 439.25
$ws.Range("K127").Value = 1317.75
$ws.Range("M127").Value = 3642.25

$ws.Range("H132").Value = 1299.5405
$ws.Range("I132").Value = 1252.3334
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 3757.0002
$ws.Range("L132").Value = 8997
$ws.Range("M132").Value = -1227.0002
$ws.Range("N132").Value = -14057

$ws.Range("H137").Value = 47623000
$ws.Range("I137").Value = 142861280
$ws.Range("J137").Value = 3859.7144
$ws.Range("K137").Value = 428583840
$ws.Range("L137").Value = 11579.1432
$ws.Range("M137").Value = -428581290
$ws.Range("N137").Value = -16679.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 14494330
$ws.Range("I74").Value = 15874552
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 15874552
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -15873678
$ws.Range("N74").Value = -3748

$ws.Range("H77").Value = 14494330
$ws.Range("I77").Value = 15874552
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 79372760
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -79368392
$ws.Range("N77").Value = -18736

$ws.Range("H122").Value = 250000960
$ws.Range("I122").Value = 1284.6666
$ws.Range("J122").Value = 1000000000
$ws.Range("K122").Value = 3853.9998
$ws.Range("L122").Value = 3000000000
$ws.Range("M122").Value = -1403.9998
$ws.Range("N122").Value = -3000004900

$ws.Range("H132").Value = 3298.675
$ws.Range("I132").Value = 2278.7144
$ws.Range("J132").Value = 10438.4
$ws.Range("K132").Value = 6836.1432
$ws.Range("L132").Value = 31315.2
$ws.Range("M132").Value = -4306.1432
$ws.Range("N132").Value = -36375.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2351.647
$ws.Range("J16").Value = 4540.5
$ws.Range("L16").Value = 4540.5
$ws.Range("N16").Value = -5114.5

$ws.Range("H86").Value = 16252
$ws.Range("J86").Value = 16252
$ws.Range("L86").Value = 16252
$ws.Range("N86").Value = -18498

$ws.Range("H89").Value = 16252
$ws.Range("J89").Value = 16252
$ws.Range("L89").Value = 81260
$ws.Range("N89").Value = -92492

$ws.Range("H113").Value = 2351.647
$ws.Range("J113").Value = 4540.5
$ws.Range("L113").Value = 4540.5
$ws.Range("N113").Value = -8880.5

$ws.Range("H134").Value = 3386.7144
$ws.Range("I134").Value = 1893.1333
$ws.Range("J134").Value = 7120.6665
$ws.Range("K134").Value = 5679.3999
$ws.Range("L134").Value = 21361.9995
$ws.Range("M134").Value = -3144.3999
$ws.Range("N134").Value = -26431.9995

$ws.Range("H141").Value = 209989.33
$ws.Range("J141").Value = 209989.33
$ws.Range("L141").Value = 209989.33
$ws.Range("N141").Value = -220349.33

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 116.666664
$ws.Range("I2").Value = 32.57143
$ws.Range("J2").Value = 234.4
$ws.Range("K2").Value = 195.42858
$ws.Range("L2").Value = 1406.4
$ws.Range("M2").Value = -82.42858000000001
$ws.Range("N2").Value = -1632.4

$ws.Range("H38").Value = 27
$ws.Range("I38").Value = 38.77778
$ws.Range("J38").Value = 13.75
$ws.Range("K38").Value = 116.33334
$ws.Range("L38").Value = 41.25
$ws.Range("M38").Value = 230.66666
$ws.Range("N38").Value = -735.25

$ws.Range("H70").Value = 10753.25
$ws.Range("I70").Value = 9333
$ws.Range("J70").Value = 15014
$ws.Range("K70").Value = 27999
$ws.Range("L70").Value = 45042
$ws.Range("M70").Value = -27684
$ws.Range("N70").Value = -45672

$ws.Range("H73").Value = 10753.25
$ws.Range("I73").Value = 9333
$ws.Range("J73").Value = 15014
$ws.Range("K73").Value = 27999
$ws.Range("L73").Value = 45042
$ws.Range("M73").Value = -26907
$ws.Range("N73").Value = -47226

$ws.Range("H132").Value = 3895.652
$ws.Range("I132").Value = 2277
$ws.Range("K132").Value = 20493
$ws.Range("M132").Value = -17963

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 13336740
$ws.Range("I102").Value = 16669675
$ws.Range("K102").Value = 16669675
$ws.Range("M102").Value = -16668053

$ws.Range("H122").Value = 3651.6667
$ws.Range("I122").Value = 1380.4
$ws.Range("K122").Value = 4141.200000000001
$ws.Range("M122").Value = -1691.200000000001

$ws.Range("H132").Value = 1317187.1
$ws.Range("I132").Value = 2191226.2
$ws.Range("J132").Value = 6128.5
$ws.Range("K132").Value = 6573678.600000001
$ws.Range("L132").Value = 18385.5
$ws.Range("M132").Value = -6571148.600000001
$ws.Range("N132").Value = -23445.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3825
$ws.Range("I16").Value = 2858.3333
$ws.Range("J16").Value = 5275
$ws.Range("K16").Value = 2858.3333
$ws.Range("L16").Value = 5275
$ws.Range("M16").Value = -2688.3333
$ws.Range("N16").Value = -5615

$ws.Range("H20").Value = 26167.5
$ws.Range("J20").Value = 26666.666
$ws.Range("L20").Value = 26666.666
$ws.Range("N20").Value = -27118.666

$ws.Range("H22").Value = 6531.2
$ws.Range("J22").Value = 7226.625
$ws.Range("L22").Value = 7226.625
$ws.Range("N22").Value = -7816.625

$ws.Range("H27").Value = 6531.2
$ws.Range("J27").Value = 7226.625
$ws.Range("L27").Value = 7226.625
$ws.Range("N27").Value = -7440.625

$ws.Range("H55").Value = 2382708.5
$ws.Range("I55").Value = 4167019.8
$ws.Range("K55").Value = 4167019.8
$ws.Range("M55").Value = -4166846.8

$ws.Range("H122").Value = 4454.4443
$ws.Range("I122").Value = 1257
$ws.Range("K122").Value = 3771
$ws.Range("M122").Value = -1321

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5758.35
$ws.Range("I122").Value = 1536.8
$ws.Range("J122").Value = 18423
$ws.Range("K122").Value = 4610.4
$ws.Range("L122").Value = 55269
$ws.Range("M122").Value = -2160.4
$ws.Range("N122").Value = -60169

$ws.Range("H132").Value = 5706.113
$ws.Range("I132").Value = 4034.7292
$ws.Range("K132").Value = 12104.1876
$ws.Range("M132").Value = -9574.187600000001

$ws.Range("H136").Value = 2240.4644
$ws.Range("I136").Value = 1653.12
$ws.Range("J136").Value = 7135
$ws.Range("K136").Value = 4959.36
$ws.Range("L136").Value = 21405
$ws.Range("M136").Value = -2409.36
$ws.Range("N136").Value = -26505
